# Add a new "Pop Quiz" (singers) block of 4 rows to the MGQ_dict sheet,
# mirroring the existing Metal/Classical/Jazz/HipHop quiz blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30: TESTNAME_POP
$ws.Range("A30").Value = "TESTNAME_POP"
$ws.Range("B30").Value = "Pop Quiz"
$ws.Range("C30").Value = "Pop Quiz"

# Row 31: INSTRUCTIONS_POP (long wrapped instructions, row height 60)
$ws.Range("A31").Value = "INSTRUCTIONS_POP"
$ws.Range("B31").Value = "Sie werden eine Liste mit Name sehen und sollen dort ankreuzen, welche davon **Sänger:innen**  sind. Sie haben dazu {{time_out}} Sekunden Zeit."
$ws.Range("C31").Value = "You will be presented with a list of names and you are asked to select all names which are **singers**. You have {{time_out}} seconds to do this."
$ws.Range("B31").WrapText = $true
$ws.Rows.Item(31).RowHeight = 60

# Row 32: PROMPT_POP
$ws.Range("A32").Value = "PROMPT_POP"
$ws.Range("B32").Value = "Bitte wählen Sie alle **Sänger:innen** aus der untenstehenden Liste aus.  Sie haben {{time_out}} Sekunden Zeit."
$ws.Range("C32").Value = "Please select all  **singers**. You have {{time_out}} seconds."

# Row 33: WELCOME_POP
$ws.Range("A33").Value = "WELCOME_POP"
$ws.Range("B33").Value = "Willkommen zum Pop Quiz"
$ws.Range("C33").Value = "Welcome to the Pop Quiz!"

# Reflect the author's final selection/viewport after the edit.
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C33").Select()
